$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. "5.20", "307.79") are stored as literal text, matching the
# original inlineStr cells instead of being coerced to numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.217.06"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "2.259.15"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "307.79"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "97.09"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "0.525"
$ws.Range("E7").Value = "  -0.91%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "34.74"
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "6.82"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "2.612.04"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "14.59"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "2.252.28"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "42.082.74"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "12.24"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "5.92"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "67.49"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "235.55"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "23.46"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "36.82"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.13"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "9.52"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "163.86"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").Value = "5.20"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "3.07"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").Value = "17.48"
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("D36").Value = "0.0727"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "1.80"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").Value = "4.14"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  -6.71%  "
$ws.Range("D43").Value = "1.940.20"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "0.0282"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "18.60"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "9.69"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").Value = "53.97"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").Value = "2.484.87"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "71.35"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "91.39"
$ws.Range("E51").Value = "  -0.69%  "

# Restore the original (General) formatting now that the text values
# are committed, so no residual text-format style remains on the cells.
$dataRange.ClearFormats()
